$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the header "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 14:52"

# 2) Update the per-country statistics that changed with this data refresh.
#    (country name stays in column A; only the numeric columns B:H move)
$ws.Range("E42").Value2 = 7174
$ws.Range("F42").Value2 = 53
$ws.Range("G42").Value2 = 1
$ws.Range("H42").Value2 = 195

$ws.Range("E52").Value2 = 2218
$ws.Range("F52").Value2 = 61
$ws.Range("G52").Value2 = 5
$ws.Range("H52").Value2 = 177

$ws.Range("E59").Value2 = 2087
$ws.Range("G59").Value2 = 4
$ws.Range("H59").Value2 = 84

$ws.Range("D64").Value2 = 602
$ws.Range("E64").Value2 = 1708

$ws.Range("B66").Value2 = 2009
$ws.Range("C66").Value2 = 28
$ws.Range("D66").Value2 = 982
$ws.Range("E66").Value2 = 976
$ws.Range("G66").Value2 = 1
$ws.Range("H66").Value2 = 51

$ws.Range("B81").Value2 = 1330
$ws.Range("C81").Value2 = 51
$ws.Range("E81").Value2 = 1109

$ws.Range("B84").Value2 = 1279
$ws.Range("C84").Value2 = 125
$ws.Range("D84").Value2 = 134
$ws.Range("E84").Value2 = 1135
$ws.Range("G84").Value2 = 1
$ws.Range("H84").Value2 = 10

$ws.Range("B114").Value2 = 414
$ws.Range("C114").Value2 = 46
$ws.Range("E114").Value2 = 300

$ws.Range("B118").Value2 = 336
$ws.Range("C118").Value2 = 16
$ws.Range("E118").Value2 = 233

$ws.Range("B148").Value2 = 90
$ws.Range("C148").Value2 = 2
$ws.Range("E148").Value2 = 25

$ws.Range("B150").Value2 = 88
$ws.Range("C150").Value2 = 6
$ws.Range("E150").Value2 = 86

$ws.Range("B160").Value2 = 82
$ws.Range("C160").Value2 = 18
$ws.Range("E160").Value2 = 70
$ws.Range("G160").Value2 = 1
$ws.Range("H160").Value2 = 2

# 3) Re-sort the country table (rows 4-216) by "Casos totales" (column B)
#    descending, same as the site does on every data refresh.
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)
